$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.830.87'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.638.21'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.258'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0639'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.81'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("D13").Value = '1.863.68'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '1.636.18'
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = '25.845.63'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.79%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0497'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.909'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = '1.133.17'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.54'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.546'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0157'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("D45").Value = '1.773.04'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  +3.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.418'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("E50").Value = '  +3.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.43%  '
